$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "2022-12-31" table column header to "#2022-12-31" (table1 / CO__crediti2022)
$ws.Range("B16").Value = "#2022-12-31"

# Update settings table values
$ws.Range("A3").Value = '$$ mod GenericMovementsH.Settings'
$ws.Range("B5").Value = "BS_CREDIT__TRADERECEIVABLECREDITS"
$ws.Range("A6").Value = "vs type"
$ws.Range("B6").Style = "Normal"
$ws.Range("B6").Value = "Bs_Cash__BankAccount_FinancialAccount"

# Drop the leftover number-format style from the "mamma"/"gino" description cells
$ws.Range("C17").Style = "Normal"
$ws.Range("C18").Style = "Normal"

# Remove row 7 ("opposite type" / "Cash") - everything below shifts up by one row
$ws.Rows("7:7").Delete()

# Restore the selection shown in the saved file
$ws.Range("A6").Select() | Out-Null
